# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 21:55"

# Murcia's case count (897) overtakes Valladolid (886), so Murcia now sorts
# ahead of Valladolid in the descending-by-"Casos totales" list.
# Row 27 becomes Murcia with its updated figures, row 28 becomes Valladolid
# with the figures that used to belong to the old row 27 (Valladolid unchanged).
$ws.Range("A27").Value = "Murcia"
$ws.Range("B27").Value = 897
$ws.Range("C27").Value = 17
$ws.Range("D27").Value = 855
$ws.Range("E27").Value = 25

$ws.Range("A28").Value = "Valladolid"
$ws.Range("B28").Value = 886
$ws.Range("C28").Value = 127
$ws.Range("D28").Value = 702
$ws.Range("E28").Value = 57
